$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.469.41"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.007.23"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.40"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.92"
$ws.Range("E7").Value = "  +5.39%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0805"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.99"
$ws.Range("E12").Value = "  +8.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.90"
$ws.Range("E13").Value = "  +3.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.856"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.297.47"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.47"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.005.28"
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.383.87"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.60"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0870"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232.34"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  +4.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.37"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.02"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("E30").Value = "  +18.75%  "
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("E32").Value = "  +4.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0628"
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("E34").Value = "  +6.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.36"
$ws.Range("E35").Value = "  +5.19%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.55"
$ws.Range("E39").Value = "  -3.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0985"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("E42").Value = "  +1.93%  "
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.88"
$ws.Range("E44").Value = "  +5.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.46"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.386.89"
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.05"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.29"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.85"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.74"
$ws.Range("E50").Value = "  +6.67%  "
$ws.Range("E51").Value = "  +12.38%  "
